# Add a new "FAVÖK" row to both the DATA_RAW sheet (row 19) and the
# "gelir tablosu (çeyreklik)" sheet (row 7). The same 60 numeric values
# (columns B..BI) are written on both sheets, with column A holding the
# new shared-string label "FAVÖK".

$wb = $excel.ActiveWorkbook

$values = @(
    0,0,0,407000000,414000000,0,491000000,-743000000,377000000,0,
    606000000,562000000,736000000,0,1339000000,1390000000,1397000000,0,1390000000,1273000000,
    1459000000,0,2562000000,2734000000,3117000000,0,3998000000,4401000000,5122000000,0,
    6398000000,5863000000,5781000000,0,4672000000,4626000000,7185000000,0,7490000000,5963000000,
    8040000000,0,13350000000,12367000000,10423000000,0,10680000000,8358000000,5664000000,0,
    9751000000,9392000000,13748000000,0,17192000000,14847000000,17350000000,21455000000,20441000000,20535000000
)

function Set-FavokRow($worksheet, $row) {
    $worksheet.Cells.Item($row, 1).Value = "FAVÖK"

    for ($i = 0; $i -lt $values.Length; $i++) {
        $worksheet.Cells.Item($row, $i + 2).Value = $values[$i]
    }
}

$wsDataRaw = $wb.Worksheets.Item("DATA_RAW")
Set-FavokRow $wsDataRaw 19

$wsGelir = $wb.Worksheets.Item("gelir tablosu (çeyreklik)")
Set-FavokRow $wsGelir 7
